# Rename the worksheet from "updated" to "Tabelle1" (input-table clean-up).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Tabelle1"
